# MitsosBarton2006Ex312 "M_Stationary generator alpha_non_zero" — refresh the
# generator's numeric outputs (new (x,y) solution point, new MIU/Lambda/Beta/
# Gamma restriction values, and new bf/BF/alpha vectors).
#
# Several of the source values are TEXT cells that merely look like numbers
# (e.g. "-1.8", "0.74") — the original file (authored by a non-Excel writer)
# stores them as shared strings, not numeric cells. Plain `Range.Value =
# "-1.8"` gets auto-coerced to a number by Excel, so we briefly force the
# range to Text format, assign, then clear the format again (ClearFormats
# keeps the already-committed text content but drops the now-unneeded style
# index) so the cell ends up exactly like the source: a text value with the
# sheet's default (unstyled) formatting.

$wb = $excel.ActiveWorkbook

function Set-TextValue {
    param($range, [string]$value)
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

# --- Restricciones_del_lider ---
$ws = $wb.Worksheets.Item("Restricciones_del_lider")
Set-TextValue $ws.Range("A2") "0.8 - x"
Set-TextValue $ws.Range("B2") "-1.8"
Set-TextValue $ws.Range("D2") "0.74"
Set-TextValue $ws.Range("A3") "-0.8 + x"
Set-TextValue $ws.Range("B3") "-0.19999999999999996"
Set-TextValue $ws.Range("D3") "0.96"

# --- Restricciones_del_follower ---
$ws = $wb.Worksheets.Item("Restricciones_del_follower")
Set-TextValue $ws.Range("A2") "0"
Set-TextValue $ws.Range("B2") "-1"
Set-TextValue $ws.Range("D2") "0.76"
Set-TextValue $ws.Range("E2") "4.4"
Set-TextValue $ws.Range("F2") "0"
Set-TextValue $ws.Range("A3") "-0.7399999999999999 + 0.3999999999999999y"
Set-TextValue $ws.Range("B3") "-0.2600000000000001"
Set-TextValue $ws.Range("D3") "0.2"
Set-TextValue $ws.Range("E3") "0"
Set-TextValue $ws.Range("F3") "3.1"

# --- Punto_modificado ---
$ws = $wb.Worksheets.Item("Punto_modificado")
Set-TextValue $ws.Range("A2") "0.8"
Set-TextValue $ws.Range("B2") "1.85"

# --- Vector_bf  (sheet index 5 — "Vector_BF" below differs only by case) ---
$ws = $wb.Worksheets.Item(5)
Set-TextValue $ws.Range("A2") "-9.78325"

# --- Vector_BF  (sheet index 6) ---
$ws = $wb.Worksheets.Item(6)
Set-TextValue $ws.Range("A2") "1.1500000000000006"
Set-TextValue $ws.Range("A3") "-49.161"

# --- Vector_Alpha (plain numeric cell, not text) ---
$ws = $wb.Worksheets.Item("Vector_Alpha")
$ws.Range("A2").Value = 0.6000000000000001
